$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.508.86"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.685.01"
$ws.Range("E3").Value = "  +5.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "419.54"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.23"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.858.92"
$ws.Range("E7").Value = "  +10.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.640"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.756"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.179"
$ws.Range("E11").Value = "  +8.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000383"
$ws.Range("E12").Value = "  +42.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.68"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.64"
$ws.Range("E14").Value = "  +7.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.260.10"
$ws.Range("E15").Value = "  +4.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.139"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.750.43"
$ws.Range("E17").Value = "  +7.14%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.48"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.87"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "66.416.89"
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.79"
$ws.Range("E22").Value = "  -4.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.61"
$ws.Range("E23").Value = "  +23.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.08"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("E25").Value = "  -5.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.34"
$ws.Range("E26").Value = "  +8.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.26"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.28"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("E29").Value = "  +4.16%  "
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("E31").Value = "  +7.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.22"
$ws.Range("E33").Value = "  -5.23%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.42"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.15"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0490"
$ws.Range("E38").Value = "  -5.11%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.98"
$ws.Range("E39").Value = "  +44.35%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  +37.72%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0748"
$ws.Range("E41").Value = "  +11.93%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.147"
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.994"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.09"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("E46").Value = "  -3.64%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  -8.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.57"
$ws.Range("E49").Value = "  -7.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.304"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.159"
$ws.Range("E51").Value = "  +12.03%  "
